$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 19; this pushes old rows 19,20 down to 20,21
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new weekly data point
$ws.Cells.Item(19, 1).Value = 1
$ws.Cells.Item(19, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(19, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(19, 4).Value = 44714
$ws.Cells.Item(19, 4).NumberFormat = $ws.Cells.Item(20, 4).NumberFormat
$ws.Cells.Item(19, 5).Value = 15
$ws.Cells.Item(19, 6).Value = 100112003
$ws.Cells.Item(19, 7).Value = "Ajo"
$ws.Cells.Item(19, 8).Value = "Chino"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 400
$ws.Cells.Item(19, 11).Value = 19000
$ws.Cells.Item(19, 12).Value = 20000
$ws.Cells.Item(19, 13).Value = 19500
$ws.Cells.Item(19, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(19, 15).Value = "China"
$ws.Cells.Item(19, 16).Value = 1950
$ws.Cells.Item(19, 17).Value = 10
$ws.Cells.Item(19, 18).Value = "Hortaliza"
